# Refresh the crypto price/volume table to the latest scrape.
#
# Most D/E cells hold plain text in the source file (t="inlineStr"),
# including values that look numeric (e.g. "479.01"). Assigning a
# bare numeric-looking string via .Value lets Excel auto-convert it
# to a real number, which would change the cell's stored type. To
# keep those cells as text we prefix the literal with a leading
# apostrophe (Excel's standard "treat as text" entry marker), then
# strip the Text number-format it applies via ClearFormats() so the
# cell's style index is left exactly as it was.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.031.41"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "2.394.92"
$ws.Range("E3").Value = "  -4.47%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'479.01"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").Value = "'148.33"
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "'0.499"
$ws.Range("E8").Value = "  -2.23%  "
$ws.Range("D9").Value = "2.396.33"
$ws.Range("E9").Value = "  -4.99%  "
$ws.Range("D10").Value = "'0.0975"
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("D11").Value = "'5.44"
$ws.Range("E11").Value = "  -4.32%  "
$ws.Range("D12").Value = "'0.322"
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("D13").Value = "'0.124"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").Value = "2.815.03"
$ws.Range("E14").Value = "  -4.11%  "
$ws.Range("D15").Value = "56.383.26"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "'20.34"
$ws.Range("E16").Value = "  -3.88%  "
$ws.Range("D17").Value = "'0.0000132"
$ws.Range("E17").Value = "  -2.85%  "
$ws.Range("D18").Value = "2.392.60"
$ws.Range("E18").Value = "  -5.07%  "
$ws.Range("D19").Value = "'4.49"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "'315.09"
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("D21").Value = "'9.77"
$ws.Range("E21").Value = "  -5.97%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "'5.68"
$ws.Range("E23").Value = "  -2.32%  "
$ws.Range("D24").Value = "'56.96"
$ws.Range("E24").Value = "  -2.96%  "
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("E26").Value = "  -4.03%  "
$ws.Range("D27").Value = "'0.158"
$ws.Range("E27").Value = "  -5.20%  "
$ws.Range("D28").Value = "2.503.56"
$ws.Range("E28").Value = "  -4.11%  "
$ws.Range("D29").Value = "'7.29"
$ws.Range("E29").Value = "  -3.47%  "
$ws.Range("D30").Value = "0.0₃0773"
$ws.Range("E30").Value = "  -2.36%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "'148.34"
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("D33").Value = "'17.97"
$ws.Range("E33").Value = "  -2.15%  "
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("D35").Value = "'5.10"
$ws.Range("E35").Value = "  -2.38%  "
$ws.Range("D36").Value = "'1.10"
$ws.Range("E36").Value = "  -3.65%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'3.58"
$ws.Range("E37").Value = "  -3.53%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'0.844"
$ws.Range("E38").Value = "  -3.07%  "
$ws.Range("D39").Value = "'33.47"
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.34"
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("D42").Value = "'3.38"
$ws.Range("E42").Value = "  -4.60%  "
$ws.Range("D43").Value = "'0.0541"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("D44").Value = "'0.0945"
$ws.Range("E44").Value = "  +3.87%  "
$ws.Range("E45").Value = "  -5.58%  "
$ws.Range("D46").Value = "'10.21"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "'4.64"
$ws.Range("E47").Value = "  -3.59%  "
$ws.Range("D48").Value = "'253.28"
$ws.Range("E48").Value = "  -4.59%  "
$ws.Range("D49").Value = "'0.0223"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("E50").Value = "  -3.81%  "
$ws.Range("D51").Value = "1.774.92"
$ws.Range("E51").Value = "  -7.82%  "

# Restore the default (General) style on the forced-text cells above.
$textForcedCells = @("D5", "D6", "D7", "D8", "D10", "D11", "D12", "D13", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D29", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49")
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).ClearFormats()
}
